# Extracted functionality to the DX12Device class.
#
# UC001 ("Engine Startup") used to spell out its Inputs/Outputs as separate
# Heading3 sub-sections. That breakdown now lives on the (soon to be
# extracted) DX12Device class, so here the Description paragraph instead
# gets a short summary sentence appended to it, and the old Inputs/Outputs
# heading+body paragraphs are removed entirely.

$d = $word.ActiveDocument

# 1. Remove the "Inputs" / "Outputs" heading+body paragraphs that followed
#    the UC001 Description paragraph.
$findRange = $d.Content
$found = $findRange.Find.Execute( `
    "Inputs`rThe width and height of the display window.`rOutputs`rAll necessary resources are initialized, the display window is created, everything is ready to run the game loop.`r", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $delRange = $d.Range($findRange.Start, $findRange.End)
    $delRange.Delete()
}

# 2. Append the new summary sentence to the end of the UC001 Description
#    paragraph text.
$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute( `
    "The main goal of this Use Case is to get the engine up and running.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $insRange = $d.Range($findRange2.End, $findRange2.End)
    $insRange.InsertAfter(" The procedure should create and display the window of the given height. It should also create and initialize the entire render system.")
}
